# GW22 league table update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in week 21 scores for the three players
$ws.Range("B19").Value = 105
$ws.Range("C19").Value = 115
$ws.Range("D19").Value = 90

# Reset the view: normal zoom and move the selection back to the top
$excel.ActiveWindow.Zoom = 100
$ws.Range("E21").Select()
